$wb = $excel.ActiveWorkbook

# Helper: write a text value to a cell without Excel's automatic
# number/string type inference kicking in (e.g. "0123456789." or "0.0"
# would otherwise silently become numeric values). We stage the text via
# a quoted formula, then collapse it to a plain static value with
# PasteSpecial (xlPasteValues = -4163) so no formula / extra number
# format survives in the saved cell.
function Set-TextValue {
    param($range, [string]$text)
    $range.Formula = '="' + $text.Replace('"', '""') + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# --- Typography sheet: fill in Wildcard Characters / Widget Wildcard
#     Characters for the "Default" typography (row 4) ---
$wsTypography = $wb.Worksheets.Item("Typography")
$wsTypography.Range("G4").Value = "."
Set-TextValue $wsTypography.Range("H4") "0123456789."

# --- Translation sheet: design the data display for VehicleInfo ---
$wsTranslation = $wb.Worksheets.Item("Translation")

# Row 4 (existing row, values replaced)
$wsTranslation.Range("B4").Value = "SingleUseId1"
$wsTranslation.Range("C4").Value = "Default"
$wsTranslation.Range("D4").Value = "Center"
$wsTranslation.Range("E4").Value = "LTR"
$wsTranslation.Range("F4").Value = "Vehicle Speed"

# Row 5 (new)
$wsTranslation.Range("B5").Value = "SingleUseId3"
$wsTranslation.Range("C5").Value = "Default"
$wsTranslation.Range("D5").Value = "Left"
$wsTranslation.Range("E5").Value = "LTR"
$wsTranslation.Range("F5").Value = "Power Percentage"

# Row 6 (new)
$wsTranslation.Range("B6").Value = "SingleUseId4"
$wsTranslation.Range("C6").Value = "Default"
$wsTranslation.Range("D6").Value = "Left"
$wsTranslation.Range("E6").Value = "LTR"
$wsTranslation.Range("F6").Value = "SOC Percentage"

# Row 7 (new)
$wsTranslation.Range("B7").Value = "SingleUseId5"
$wsTranslation.Range("C7").Value = "Default"
$wsTranslation.Range("D7").Value = "Left"
$wsTranslation.Range("E7").Value = "LTR"
$wsTranslation.Range("F7").Value = "Vehicle Status <d>"

# Row 8 (new)
$wsTranslation.Range("B8").Value = "SingleUseId6"
$wsTranslation.Range("C8").Value = "Default"
$wsTranslation.Range("D8").Value = "Left"
$wsTranslation.Range("E8").Value = "LTR"
$wsTranslation.Range("F8").Value = "VCU Errors <value>"

# Row 9 (new)
$wsTranslation.Range("B9").Value = "SingleUseId7"
$wsTranslation.Range("C9").Value = "Default"
$wsTranslation.Range("D9").Value = "Center"
$wsTranslation.Range("E9").Value = "LTR"
$wsTranslation.Range("F9").Value = "<d> Km/h"

# Row 10 (new)
$wsTranslation.Range("B10").Value = "SingleUseId8"
$wsTranslation.Range("C10").Value = "Default"
$wsTranslation.Range("D10").Value = "Left"
$wsTranslation.Range("E10").Value = "LTR"
Set-TextValue $wsTranslation.Range("F10") "0.0"

# Row 11 (new)
$wsTranslation.Range("B11").Value = "SingleUseId9"
$wsTranslation.Range("C11").Value = "Default"
$wsTranslation.Range("D11").Value = "Left"
$wsTranslation.Range("E11").Value = "LTR"
Set-TextValue $wsTranslation.Range("F11") "0"

# Row 12 (new)
$wsTranslation.Range("B12").Value = "SingleUseId10"
$wsTranslation.Range("C12").Value = "Default"
$wsTranslation.Range("D12").Value = "Left"
$wsTranslation.Range("E12").Value = "LTR"
Set-TextValue $wsTranslation.Range("F12") "0"

# Row 13 (new)
$wsTranslation.Range("B13").Value = "SingleUseId11"
$wsTranslation.Range("C13").Value = "Default"
$wsTranslation.Range("D13").Value = "Center"
$wsTranslation.Range("E13").Value = "LTR"
$wsTranslation.Range("F13").Value = "<> %"

$excel.CutCopyMode = $false
